# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet - this shifts the old N/O/P ("Late" / "Heading" / "Outstanding")
# columns one place to the right (-> O/P/Q) and leaves the freshly
# inserted column N empty. Also switches the active sheet/tab from
# "Transactions" back to "Repayment Schedule".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N (14th column); existing N..P data
# shifts right to O..Q.
$ws.Columns("N").Insert()

# Give the newly inserted column roughly the same width as its neighbour
# (column M) rather than leaving it at the sheet default.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab again (was "Transactions").
$ws.Activate()
$ws.Range("R9").Select()
